$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, pushing existing rows 38-47 down to 39-48.
$ws.Rows.Item(38).Insert()

# Populate the new row 38 with the new weekly data point.
# The constant columns (market/region/category/unit/origin/classification)
# mirror the surrounding rows for this series.
$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value = "Ñuble"
$ws.Cells.Item(38, 4).Value = 44694
$ws.Cells.Item(38, 5).Value = 16
$ws.Cells.Item(38, 6).Value = 100112040
$ws.Cells.Item(38, 7).Value = "Cilantro"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 120
$ws.Cells.Item(38, 11).Value = 550
$ws.Cells.Item(38, 12).Value = 600
$ws.Cells.Item(38, 13).Value = 575
$ws.Cells.Item(38, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(38, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(38, 16).Value = 575
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"
